# "draw functions and map" -- rework the Enemy sheet's size table and
# rename the "map" sheet's enemy-list column to enemy-id.

$wb = $excel.ActiveWorkbook

# --- map sheet: the 3rd column used to be a list of enemy ids
#     ("enemy_list"); it is now a single enemy id ("enemy_id"). ---
$mapSheet = $wb.Worksheets.Item("map")
$mapSheet.Range("C3").Value = "enemy_id"

# --- Enemy sheet: add three more enemy "size" presets (rows 5-7),
#     following on from the existing row 4 (id=1, size="10,10"). ---
$enemySheet = $wb.Worksheets.Item("Enemy")

$newRows = @(
    @(2, 1, 5, 1, "10,11"),
    @(3, 1, 5, 1, "10,12"),
    @(4, 1, 5, 1, "10,13")
)

$r = 5
foreach ($row in $newRows) {
    $enemySheet.Cells.Item($r, 1).Value = $row[0]
    $enemySheet.Cells.Item($r, 2).Value = $row[1]
    $enemySheet.Cells.Item($r, 3).Value = $row[2]
    $enemySheet.Cells.Item($r, 4).Value = $row[3]
    $enemySheet.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# --- The Enemy sheet becomes the active/selected sheet (was "map"). ---
$enemySheet.Select()
